$d = $word.ActiveDocument

# Paragraph 3 ("O metodo agil pode ser aplicado...") currently ends with the
# _GoBack bookmark right after its single run. Replace that whole paragraph
# (text + paragraph mark) with three paragraphs:
#   1) the same text, now ending its own paragraph;
#   2) a new blank paragraph;
#   3) a new paragraph with the additional narrative (two runs), followed by
#      the original _GoBack bookmark, followed by a closing ". " run.
$p3 = $d.Paragraphs.Item(3)
$target = $d.Range($p3.Range.Start, $p3.Range.End)

$xml = '<w:p w14:paraId="224E1991" w14:textId="279A0AFD" w:rsidR="00D52F9E" w:rsidRDefault="0004709D" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>O método ágil pode ser aplicado ao projeto, por conta de ter um ciclo de vida evolucionário, assim gerando uma boa prática em gestão do projeto, e ajuda na parte da empresa ter a necessidade de diminuir custos, otimizando recursos e tempo, utilizando uma equipe pequena e com múltiplas funções com foco a atingir o objetivo do projeto.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Como foi aplicado um método ágil ao projeto a equipe tende a ser menor para diminuição de custos, e tende a ser multifuncionais. Tendo como base o dono da empresa exercendo o papel de PO, pois foi com ele que teve a entrevista para ter ideias iniciais e entendimento dos problemas que a empresa enfrenta, já na parte do time vai conter um ScrumMaster para ter uma segurança, sendo assim, impedindo conflitos externos dentro do projeto e garantindo ao máximo possível os prazos designados, e por fim o time fullstack onde vão conter 3 full stacks para que seja possível criar os protótipos o mais rápido possível e garantir os testes, onde vão liberar a fase de implantação desses pequenos softwares que vão ser implementados durante o ciclo do projeto, implementando a resolução de cada problema que foi colocado pelo PO</w:t></w:r><w:r><w:t xml:space="preserve"> até atingir um produto final totalmente testado, com algumas partes já implementadas assim agilizando o processo de treinamento</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p>'
$target.InsertXML($xml)

Write-Host "done"
